$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 2.242386666666667
$ws.Cells.Item(2, 8).Value = 6.72716
$ws.Cells.Item(2, 9).Value = 0.04442500453715972
$ws.Cells.Item(2, 10).Value = 0.04442500453715972
$ws.Cells.Item(2, 13).Value = 13.713764
$ws.Cells.Item(2, 14).Value = 41.141292
$ws.Cells.Item(2, 15).Value = 0.0901423721847377
$ws.Cells.Item(2, 16).Value = 0.0901423721847377
$ws.Cells.Item(2, 17).Value = 30.75156154341333
$ws.Cells.Item(2, 18).Value = 276.76405389072
$ws.Cells.Item(2, 19).Value = 0.004004575293297312
$ws.Cells.Item(2, 20).Value = 0.004004575293297312
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 2.242386666666667
$ws.Cells.Item(3, 8).Value = 6.72716
$ws.Cells.Item(3, 9).Value = 0.04442500453715972
$ws.Cells.Item(3, 10).Value = 0.04442500453715972
$ws.Cells.Item(3, 14).Value = 84.55600199999999
$ws.Cells.Item(3, 15).Value = 0.1852659027513629
$ws.Cells.Item(3, 16).Value = 0.1852659027513629
$ws.Cells.Item(3, 17).Value = 63.20241715714666
$ws.Cells.Item(3, 18).Value = 568.8217544143199
$ws.Cells.Item(3, 19).Value = 0.008230438570310288
$ws.Cells.Item(3, 20).Value = 0.008230438570310288
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 2.242386666666667
$ws.Cells.Item(4, 8).Value = 6.72716
$ws.Cells.Item(4, 9).Value = 0.04442500453715972
$ws.Cells.Item(4, 10).Value = 0.04442500453715972
$ws.Cells.Item(4, 13).Value = 21.07704566666666
$ws.Cells.Item(4, 14).Value = 63.23113699999999
$ws.Cells.Item(4, 15).Value = 0.1385421898057586
$ws.Cells.Item(4, 16).Value = 0.1385421898057586
$ws.Cells.Item(4, 17).Value = 47.26288617565777
$ws.Cells.Item(4, 18).Value = 425.3659755809199
$ws.Cells.Item(4, 19).Value = 0.00615473741070887
$ws.Cells.Item(4, 20).Value = 0.006154737410708869
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 2.242386666666667
$ws.Cells.Item(5, 8).Value = 6.72716
$ws.Cells.Item(5, 9).Value = 0.04442500453715972
$ws.Cells.Item(5, 10).Value = 0.04442500453715972
$ws.Cells.Item(5, 13).Value = 89.15834833333334
$ws.Cells.Item(5, 14).Value = 267.475045
$ws.Cells.Item(5, 15).Value = 0.5860495352581409
$ws.Cells.Item(5, 16).Value = 0.5860495352581407
$ws.Cells.Item(5, 17).Value = 199.9274915246889
$ws.Cells.Item(5, 18).Value = 1799.3474237222
$ws.Cells.Item(5, 19).Value = 0.02603525326284326
$ws.Cells.Item(5, 20).Value = 0.02603525326284324
$ws.Cells.Item(6, 9).Value = 0.4052409520727612
$ws.Cells.Item(6, 10).Value = 0.4052409520727612
$ws.Cells.Item(6, 13).Value = 13.713764
$ws.Cells.Item(6, 14).Value = 41.141292
$ws.Cells.Item(6, 15).Value = 0.0901423721847377
$ws.Cells.Item(6, 16).Value = 0.0901423721847377
$ws.Cells.Item(6, 17).Value = 280.5130175541826
$ws.Cells.Item(6, 18).Value = 2524.617157987644
$ws.Cells.Item(6, 19).Value = 0.03652938072624029
$ws.Cells.Item(6, 20).Value = 0.03652938072624029
$ws.Cells.Item(7, 9).Value = 0.4052409520727612
$ws.Cells.Item(7, 10).Value = 0.4052409520727612
$ws.Cells.Item(7, 14).Value = 84.55600199999999
$ws.Cells.Item(7, 15).Value = 0.1852659027513629
$ws.Cells.Item(7, 16).Value = 0.1852659027513629
$ws.Cells.Item(7, 18).Value = 5188.741604421113
$ws.Cells.Item(7, 19).Value = 0.07507733081758189
$ws.Cells.Item(7, 20).Value = 0.07507733081758189
$ws.Cells.Item(8, 9).Value = 0.4052409520727612
$ws.Cells.Item(8, 10).Value = 0.4052409520727612
$ws.Cells.Item(8, 13).Value = 21.07704566666666
$ws.Cells.Item(8, 14).Value = 63.23113699999999
$ws.Cells.Item(8, 15).Value = 0.1385421898057586
$ws.Cells.Item(8, 16).Value = 0.1385421898057586
$ws.Cells.Item(8, 17).Value = 431.1278567345898
$ws.Cells.Item(8, 18).Value = 3880.150710611308
$ws.Cells.Item(8, 19).Value = 0.05614296889913081
$ws.Cells.Item(8, 20).Value = 0.05614296889913081
$ws.Cells.Item(9, 9).Value = 0.4052409520727612
$ws.Cells.Item(9, 10).Value = 0.4052409520727612
$ws.Cells.Item(9, 13).Value = 89.15834833333334
$ws.Cells.Item(9, 14).Value = 267.475045
$ws.Cells.Item(9, 15).Value = 0.5860495352581409
$ws.Cells.Item(9, 16).Value = 0.5860495352581407
$ws.Cells.Item(9, 17).Value = 1823.720849442229
$ws.Cells.Item(9, 18).Value = 16413.48764498006
$ws.Cells.Item(9, 19).Value = 0.2374912716298082
$ws.Cells.Item(9, 20).Value = 0.2374912716298082
$ws.Cells.Item(10, 7).Value = 27.778539
$ws.Cells.Item(10, 8).Value = 83.335617
$ws.Cells.Item(10, 9).Value = 0.5503340433900792
$ws.Cells.Item(10, 10).Value = 0.5503340433900791
$ws.Cells.Item(10, 13).Value = 13.713764
$ws.Cells.Item(10, 14).Value = 41.141292
$ws.Cells.Item(10, 15).Value = 0.0901423721847377
$ws.Cells.Item(10, 16).Value = 0.0901423721847377
$ws.Cells.Item(10, 17).Value = 380.948328110796
$ws.Cells.Item(10, 18).Value = 3428.534952997164
$ws.Cells.Item(10, 19).Value = 0.0496084161652001
$ws.Cells.Item(10, 20).Value = 0.04960841616520009
$ws.Cells.Item(11, 7).Value = 27.778539
$ws.Cells.Item(11, 8).Value = 83.335617
$ws.Cells.Item(11, 9).Value = 0.5503340433900792
$ws.Cells.Item(11, 10).Value = 0.5503340433900791
$ws.Cells.Item(11, 14).Value = 84.55600199999999
$ws.Cells.Item(11, 15).Value = 0.1852659027513629
$ws.Cells.Item(11, 16).Value = 0.1852659027513629
$ws.Cells.Item(11, 17).Value = 782.9473997470259
$ws.Cells.Item(11, 18).Value = 7046.526597723233
$ws.Cells.Item(11, 19).Value = 0.1019581333634707
$ws.Cells.Item(11, 20).Value = 0.1019581333634707
$ws.Cells.Item(12, 7).Value = 27.778539
$ws.Cells.Item(12, 8).Value = 83.335617
$ws.Cells.Item(12, 9).Value = 0.5503340433900792
$ws.Cells.Item(12, 10).Value = 0.5503340433900791
$ws.Cells.Item(12, 13).Value = 21.07704566666666
$ws.Cells.Item(12, 14).Value = 63.23113699999999
$ws.Cells.Item(12, 15).Value = 0.1385421898057586
$ws.Cells.Item(12, 16).Value = 0.1385421898057586
$ws.Cells.Item(12, 17).Value = 585.4895350562808
$ws.Cells.Item(12, 18).Value = 5269.405815506528
$ws.Cells.Item(12, 19).Value = 0.07624448349591895
$ws.Cells.Item(12, 20).Value = 0.07624448349591893
$ws.Cells.Item(13, 7).Value = 27.778539
$ws.Cells.Item(13, 8).Value = 83.335617
$ws.Cells.Item(13, 9).Value = 0.5503340433900792
$ws.Cells.Item(13, 10).Value = 0.5503340433900791
$ws.Cells.Item(13, 13).Value = 89.15834833333334
$ws.Cells.Item(13, 14).Value = 267.475045
$ws.Cells.Item(13, 15).Value = 0.5860495352581409
$ws.Cells.Item(13, 16).Value = 0.5860495352581407
$ws.Cells.Item(13, 17).Value = 2476.688656353085
$ws.Cells.Item(13, 18).Value = 22290.19790717777
$ws.Cells.Item(13, 19).Value = 0.3225230103654894
$ws.Cells.Item(13, 20).Value = 0.3225230103654893
